$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "Ширина мм лице\тил"
$ws.Range("F3").Value = "Висота мм лице\тил"

$ws.Rows("4:5").Delete()

$nl = [char]10
$a7 = "Всього за дверні блоки: 0.0 грн" + $nl + "`t`tЗнижка:0.0%" + $nl + "`t`tУсього, з урахуванням знижки: 0.0 грн" + $nl + "`t`tДоставка на склад (об'єкт) без вивантаження та занесення на поверх: 0.0 грн" + $nl + "`t`tМонтаж:0.0грнЗаміри: 0.0 грн" + $nl + "`t`tВсього за послуги: 0.0 грн" + $nl + "`t`tВсього сума замовлення: 0.0 грн" + $nl + "`t`tПередплата: 0.0 % " + $nl + "`t`tПередплата: 0.0 грн " + $nl + "`t`tЗалишок: 0.0 грн"
$ws.Range("A7").Value = $a7
